$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") holds the strikeout count per outing; regen values
# to reflect strikeouts recorded as "K" rather than the old "Strike#" values.
$kValues = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
